$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G13").Value = $ws.Range("G13").Value2 + 2
$ws.Range("G16").Value = $ws.Range("G16").Value2 + 2
$ws.Range("G31").Value = $ws.Range("G31").Value2 + 2
$ws.Range("G38").Value = $ws.Range("G38").Value2 + 2
$ws.Range("G46").Value = $ws.Range("G46").Value2 + 2
$ws.Range("G79").Value = $ws.Range("G79").Value2 + 2
$ws.Range("G88").Value = $ws.Range("G88").Value2 + 2
$ws.Range("G107").Value = $ws.Range("G107").Value2 + 2
$ws.Range("G110").Value = $ws.Range("G110").Value2 + 2
$ws.Range("G115").Value = $ws.Range("G115").Value2 + 2
$ws.Range("G122").Value = $ws.Range("G122").Value2 + 2
$ws.Range("G126").Value = $ws.Range("G126").Value2 + 2
$ws.Range("G136").Value = $ws.Range("G136").Value2 + 2
$ws.Range("G152").Value = $ws.Range("G152").Value2 + 2
$ws.Range("G160").Value = $ws.Range("G160").Value2 + 2
$ws.Range("G176").Value = $ws.Range("G176").Value2 + 2
$ws.Range("G186").Value = $ws.Range("G186").Value2 + 2
$ws.Range("G207").Value = $ws.Range("G207").Value2 + 2
$ws.Range("G211").Value = $ws.Range("G211").Value2 + 2
$ws.Range("G212").Value = $ws.Range("G212").Value2 + 2
$ws.Range("G213").Value = $ws.Range("G213").Value2 + 2
$ws.Range("G214").Value = $ws.Range("G214").Value2 + 2
$ws.Range("G222").Value = $ws.Range("G222").Value2 + 2
$ws.Range("G238").Value = $ws.Range("G238").Value2 + 2
$ws.Range("G248").Value = $ws.Range("G248").Value2 + 2
$ws.Range("G250").Value = $ws.Range("G250").Value2 + 2
$ws.Range("G261").Value = $ws.Range("G261").Value2 + 2
$ws.Range("G262").Value = $ws.Range("G262").Value2 + 2
$ws.Range("G264").Value = $ws.Range("G264").Value2 + 2
$ws.Range("G273").Value = $ws.Range("G273").Value2 + 2
$ws.Range("G279").Value = $ws.Range("G279").Value2 + 2
$ws.Range("G284").Value = $ws.Range("G284").Value2 + 2
$ws.Range("G286").Value = $ws.Range("G286").Value2 + 2
$ws.Range("G287").Value = $ws.Range("G287").Value2 + 2
$ws.Range("G297").Value = $ws.Range("G297").Value2 + 2
$ws.Range("G299").Value = $ws.Range("G299").Value2 + 2
$ws.Range("G302").Value = $ws.Range("G302").Value2 + 2
$ws.Range("G308").Value = $ws.Range("G308").Value2 + 2
$ws.Range("G309").Value = $ws.Range("G309").Value2 + 2
$ws.Range("G314").Value = $ws.Range("G314").Value2 + 2
$ws.Range("G315").Value = $ws.Range("G315").Value2 + 2
$ws.Range("G319").Value = $ws.Range("G319").Value2 + 2
$ws.Range("G333").Value = $ws.Range("G333").Value2 + 2
$ws.Range("G335").Value = $ws.Range("G335").Value2 + 2
$ws.Range("G350").Value = $ws.Range("G350").Value2 + 2
